$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.304.51"
$ws.Range("E2").Value = "  -1.58%  "
$ws.Range("D3").Value = "3.069.89"
$ws.Range("E3").Value = "  -2.60%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.27"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.95"
$ws.Range("E6").Value = "  +4.74%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.540"
$ws.Range("E8").Value = "  +1.80%  "
$ws.Range("D9").Value = "3.067.81"
$ws.Range("E9").Value = "  -2.41%  "
$ws.Range("E10").Value = "  -2.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.77"
$ws.Range("E11").Value = "  -3.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.454"
$ws.Range("E12").Value = "  -1.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.20"
$ws.Range("E13").Value = "  -0.70%  "
$ws.Range("E14").Value = "  -3.59%  "
$ws.Range("D15").Value = "3.583.97"
$ws.Range("E15").Value = "  -2.42%  "
$ws.Range("E16").Value = "  -1.99%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "63.381.90"
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.12"
$ws.Range("E18").Value = "  -2.65%  "
$ws.Range("D19").Value = "3.072.17"
$ws.Range("E19").Value = "  -2.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "473.65"
$ws.Range("E20").Value = "  +0.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.50"
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.710"
$ws.Range("E22").Value = "  -3.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.48"
$ws.Range("E23").Value = "  -2.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.38"
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.87"
$ws.Range("E25").Value = "  -1.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "80.77"
$ws.Range("E26").Value = "  -0.90%  "
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.93"
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.30"
$ws.Range("E29").Value = "  -2.18%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.66"
$ws.Range("E31").Value = "  -2.30%  "
$ws.Range("E32").Value = "  -3.49%  "
$ws.Range("E33").Value = "  +3.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.06"
$ws.Range("E34").Value = "  -2.29%  "
$ws.Range("D35").Value = "0.0$([char]0x2083)0837"
$ws.Range("E35").Value = "  -1.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.04"
$ws.Range("E36").Value = "  -2.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.35"
$ws.Range("E37").Value = "  +4.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.03"
$ws.Range("E38").Value = "  -2.44%  "
$ws.Range("E39").Value = "  -5.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.58"
$ws.Range("E40").Value = "  -1.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.19"
$ws.Range("E41").Value = "  -1.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "440.99"
$ws.Range("E42").Value = "  -3.42%  "
$ws.Range("E43").Value = "  -3.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0360"
$ws.Range("E44").Value = "  -3.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "39.86"
$ws.Range("E45").Value = "  -1.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.110"
$ws.Range("E46").Value = "  +1.86%  "
$ws.Range("D47").Value = "2.787.27"
$ws.Range("E47").Value = "  -4.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.49"
$ws.Range("E48").Value = "  -2.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.49"
$ws.Range("E49").Value = "  +4.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.24"
$ws.Range("E51").Value = "  -0.30%  "
